$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.789.41"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "3.187.13"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.27"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.81"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.00%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("D12").Value = "3.737.05"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000174"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.96"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "59.813.86"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").Value = "3.194.38"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.20"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.24"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.50"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +4.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.23"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.88"
$ws.Range("E26").Value = "  +16.57%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "0.0₃0902"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.91"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.41"
$ws.Range("E30").Value = "  +4.31%  "
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  +4.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.81"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "2.765.50"
$ws.Range("E37").Value = "  +7.86%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.68"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0711"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("E40").Value = "  +3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.27"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +7.43%  "
$ws.Range("D45").Value = "3.232.14"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  +5.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.54"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.778"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.06%  "
